$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 710.5454999999999
$ws.Range("I28").Value = 668.3333
$ws.Range("J28").Value = 761.2
$ws.Range("K28").Value = 668.3333
$ws.Range("L28").Value = 761.2
$ws.Range("M28").Value = -183.3333
$ws.Range("N28").Value = -1731.2

$ws.Range("H98").Value = 35180.74
$ws.Range("I98").Value = 845.3333
$ws.Range("J98").Value = 99559.625
$ws.Range("K98").Value = 845.3333
$ws.Range("L98").Value = 99559.625
$ws.Range("M98").Value = 652.6667
$ws.Range("N98").Value = -102555.625

$ws.Range("H107").Value = 7986
$ws.Range("I107").Value = 10736.4
$ws.Range("J107").Value = 1110
$ws.Range("K107").Value = 10736.4
$ws.Range("L107").Value = 1110
$ws.Range("M107").Value = -8816.4
$ws.Range("N107").Value = -4950

$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").ClearContents()

$ws.Range("H115").Value = 2274.75
$ws.Range("I115").Value = 2274.75
$ws.Range("K115").Value = 6824.25
$ws.Range("M115").Value = -5257.25

$ws.Range("H116").Value = 2880.2693
$ws.Range("I116").Value = 2139.1333
$ws.Range("J116").Value = 3890.9092
$ws.Range("K116").Value = 2139.1333
$ws.Range("L116").Value = 3890.9092
$ws.Range("M116").Value = 1302.8667
$ws.Range("N116").Value = -10774.9092

$ws.Range("H122").Value = 35180.74
$ws.Range("I122").Value = 845.3333
$ws.Range("J122").Value = 99559.625
$ws.Range("K122").Value = 2535.9999
$ws.Range("L122").Value = 298678.875
$ws.Range("M122").Value = -85.9998999999998
$ws.Range("N122").Value = -303578.875

$ws.Range("H125").Value = 1678.5714
$ws.Range("I125").Value = 1533.3334
$ws.Range("J125").Value = 1718.1818
$ws.Range("K125").Value = 13800.0006
$ws.Range("L125").Value = 15463.6362
$ws.Range("M125").Value = -11340.0006
$ws.Range("N125").Value = -20383.6362

$ws.Range("H129").Value = 1040.6804
$ws.Range("I129").Value = 1535.8182
$ws.Range("K129").Value = 4607.4546
$ws.Range("M129").Value = 392.5454

$ws.Range("H132").Value = 21170.531
$ws.Range("I132").Value = 2676.95
$ws.Range("J132").Value = 126848.14
$ws.Range("K132").Value = 8030.849999999999
$ws.Range("L132").Value = 380544.42
$ws.Range("M132").Value = -5500.849999999999
$ws.Range("N132").Value = -385604.42

$ws.Range("H137").Value = 2963796.5
$ws.Range("I137").Value = 12821572
$ws.Range("K137").Value = 38464716
$ws.Range("M137").Value = -38462166

$ws.Range("H141").Value = 2765.5264
$ws.Range("I141").Value = 1393.8462
$ws.Range("K141").Value = 4181.5386
$ws.Range("M141").Value = 998.4614000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1366.6875
$ws.Range("I2").Value = 1373.1666
$ws.Range("J2").Value = 1347.25
$ws.Range("K2").Value = 1373.1666
$ws.Range("L2").Value = 1347.25
$ws.Range("M2").Value = -1260.1666
$ws.Range("N2").Value = -1573.25

$ws.Range("H32").Value = 3115.4775
$ws.Range("I32").Value = 1480.0469
$ws.Range("K32").Value = 1480.0469
$ws.Range("M32").Value = -1193.0469

$ws.Range("H40").Value = 29258
$ws.Range("J40").Value = 29258
$ws.Range("L40").Value = 29258
$ws.Range("N40").Value = -29610

$ws.Range("H42").Value = 20080
$ws.Range("J42").Value = 20080
$ws.Range("L42").Value = 20080
$ws.Range("N42").Value = -21052

$ws.Range("H110").Value = 1830
$ws.Range("I110").Value = 1721.5454
$ws.Range("J110").Value = 2128.25
$ws.Range("K110").Value = 1721.5454
$ws.Range("L110").Value = 2128.25
$ws.Range("M110").Value = 323.4546
$ws.Range("N110").Value = -6218.25

$ws.Range("H116").Value = 1366.6875
$ws.Range("I116").Value = 1373.1666
$ws.Range("J116").Value = 1347.25
$ws.Range("K116").Value = 1373.1666
$ws.Range("L116").Value = 1347.25
$ws.Range("M116").Value = 920.8334
$ws.Range("N116").Value = -5935.25

$ws.Range("H132").Value = 17243006
$ws.Range("I132").Value = 20834298
$ws.Range("K132").Value = 62502894
$ws.Range("M132").Value = -62500364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1366.6875
$ws.Range("I3").Value = 1373.1666
$ws.Range("J3").Value = 1347.25
$ws.Range("K3").Value = 1373.1666
$ws.Range("L3").Value = 1347.25
$ws.Range("M3").Value = -1259.1666
$ws.Range("N3").Value = -1575.25

$ws.Range("H107").Value = 1865.5555
$ws.Range("I107").Value = 1782.9
$ws.Range("J107").Value = 1968.875
$ws.Range("K107").Value = 1782.9
$ws.Range("L107").Value = 1968.875
$ws.Range("M107").Value = 137.0999999999999
$ws.Range("N107").Value = -5808.875

$ws.Range("H134").Value = 2190.0125
$ws.Range("I134").Value = 1310.8628
$ws.Range("J134").Value = 3684.5667
$ws.Range("K134").Value = 3932.588400000001
$ws.Range("L134").Value = 11053.7001
$ws.Range("M134").Value = -1397.588400000001
$ws.Range("N134").Value = -16123.7001

$ws.Range("H138").Value = 43258.402
$ws.Range("J138").Value = 43258.402
$ws.Range("L138").Value = 43258.402
$ws.Range("N138").Value = -53538.402

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1165.3478
$ws.Range("I16").Value = 1094.8948
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1094.8948
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -807.8948
$ws.Range("N16").Value = -2074

$ws.Range("H107").Value = 587.19446
$ws.Range("I107").Value = 462.9655
$ws.Range("J107").Value = 1101.8572
$ws.Range("K107").Value = 462.9655
$ws.Range("L107").Value = 1101.8572
$ws.Range("M107").Value = 1457.0345
$ws.Range("N107").Value = -4941.8572

$ws.Range("H113").Value = 1165.3478
$ws.Range("I113").Value = 1094.8948
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1094.8948
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1075.1052
$ws.Range("N113").Value = -5840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 2983.7778
$ws.Range("I119").Value = 1633.5
$ws.Range("J119").Value = 4671.625
$ws.Range("K119").Value = 4900.5
$ws.Range("L119").Value = 14014.875
$ws.Range("M119").Value = -62.5
$ws.Range("N119").Value = -23690.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1721.1
$ws.Range("I113").Value = 1785.1666
$ws.Range("J113").Value = 1625
$ws.Range("K113").Value = 1785.1666
$ws.Range("L113").Value = 1625
$ws.Range("M113").Value = 384.8334
$ws.Range("N113").Value = -5965

$ws.Range("H122").Value = 1427.88
$ws.Range("I122").Value = 1466.0555
$ws.Range("J122").Value = 1329.7142
$ws.Range("K122").Value = 4398.166499999999
$ws.Range("L122").Value = 3989.1426
$ws.Range("M122").Value = -1948.166499999999
$ws.Range("N122").Value = -8889.142599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 3000
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

$ws.Range("H40").Value = 3676.182
$ws.Range("I40").Value = 2791.75
$ws.Range("J40").Value = 4737.5
$ws.Range("K40").Value = 2791.75
$ws.Range("L40").Value = 4737.5
$ws.Range("M40").Value = -2655.75
$ws.Range("N40").Value = -5009.5

$ws.Range("H61").Value = 2906.9375
$ws.Range("I61").Value = 3066
$ws.Range("K61").Value = 3066
$ws.Range("M61").Value = -2864

$ws.Range("H113").Value = 2906.9375
$ws.Range("I113").Value = 3066
$ws.Range("K113").Value = 3066
$ws.Range("M113").Value = -896

$ws.Range("H122").Value = 1837.6666
$ws.Range("I122").Value = 1837.6666
$ws.Range("K122").Value = 5512.9998
$ws.Range("M122").Value = -3062.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 25000
$ws.Range("J39").Value = 25000
$ws.Range("L39").Value = 25000
$ws.Range("N39").Value = -25826

$ws.Range("H122").Value = 1429313.9
$ws.Range("I122").Value = 2198439.8
$ws.Range("J122").Value = 937.1429000000001
$ws.Range("K122").Value = 6595319.399999999
$ws.Range("L122").Value = 2811.4287
$ws.Range("M122").Value = -6592869.399999999
$ws.Range("N122").Value = -7711.4287

$ws.Range("H132").Value = 1012540.06
$ws.Range("I132").Value = 1243240.8
$ws.Range("J132").Value = 3224.5
$ws.Range("K132").Value = 3729722.4
$ws.Range("L132").Value = 9673.5
$ws.Range("M132").Value = -3727192.4
$ws.Range("N132").Value = -14733.5
